# Weekly data refresh: insert two new price rows (row 612/613) at the top of the
# "Naranja" table on the active sheet, shifting the existing rows 612:633 down to
# 614:635 (dimension grows from A1:T633 to A1:T635).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 612; everything below (old 612:633)
# shifts down to 614:635, carrying its formatting/values with it.
$ws.Range("A612:T613").EntireRow.Insert()

# Row 612 - new "Especial" grade entry dated 2023-05-29 (serial 45075)
$ws.Cells.Item(612, 1).Value = 7
$ws.Cells.Item(612, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(612, 3).Value = "Ñuble"
$ws.Cells.Item(612, 4).Value = 45075
$ws.Cells.Item(612, 5).Value = 16
$ws.Cells.Item(612, 6).Value = "Fruta"
$ws.Cells.Item(612, 7).Value = 100102
$ws.Cells.Item(612, 8).Value = "Cítricos"
$ws.Cells.Item(612, 9).Value = 100102005
$ws.Cells.Item(612, 10).Value = "Naranja"
$ws.Cells.Item(612, 11).Value = "Valencia"
$ws.Cells.Item(612, 12).Value = "Especial"
$ws.Cells.Item(612, 13).Value = 60
$ws.Cells.Item(612, 14).Value = 12000
$ws.Cells.Item(612, 15).Value = 12000
$ws.Cells.Item(612, 16).Value = 12000
$ws.Cells.Item(612, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(612, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(612, 19).Value = 800
$ws.Cells.Item(612, 20).Value = 15

# Row 613 - new "Primera" grade entry, same date
$ws.Cells.Item(613, 1).Value = 7
$ws.Cells.Item(613, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(613, 3).Value = "Ñuble"
$ws.Cells.Item(613, 4).Value = 45075
$ws.Cells.Item(613, 5).Value = 16
$ws.Cells.Item(613, 6).Value = "Fruta"
$ws.Cells.Item(613, 7).Value = 100102
$ws.Cells.Item(613, 8).Value = "Cítricos"
$ws.Cells.Item(613, 9).Value = 100102005
$ws.Cells.Item(613, 10).Value = "Naranja"
$ws.Cells.Item(613, 11).Value = "Valencia"
$ws.Cells.Item(613, 12).Value = "Primera"
$ws.Cells.Item(613, 13).Value = 50
$ws.Cells.Item(613, 14).Value = 10000
$ws.Cells.Item(613, 15).Value = 10000
$ws.Cells.Item(613, 16).Value = 10000
$ws.Cells.Item(613, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(613, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(613, 19).Value = 667
$ws.Cells.Item(613, 20).Value = 15
